# Applies the Aufgabe_Excel_04 corrections described in the commit:
# "Excel engl. 01 and 02 added, corrections after execution of week 1 --all"
#
# Strategy: each change is expressed as a Find/Replace over $d.Content using
# the classic positional Find.Execute(...) signature (FindText, MatchCase,
# MatchWholeWord, MatchWildcards, MatchSoundsLike, MatchAllWordForms,
# Forward, Wrap, Format, ReplaceWith, Replace) with Wrap=1 (wdFindContinue)
# and Replace=2 (wdReplaceAll). This naturally merges any runs spanned by
# the match (dropping now-pointless <w:proofErr/> gramStart/gramEnd pairs
# the same way the reference diff does).

$d = $word.ActiveDocument

function Replace-Text($find, $replace) {
    $ok = $d.Content.Find.Execute($find, $false, $false, $false, $false, $false, $true, 1, $false, $replace, 2)
    if (-not $ok) {
        throw ("Find/Replace failed for: " + $find)
    }
}

function Replace-WholeWordText($find, $replace) {
    $ok = $d.Content.Find.Execute($find, $false, $true, $false, $false, $false, $true, 1, $false, $replace, 2)
    if (-not $ok) {
        throw ("Find/Replace failed for: " + $find)
    }
}

# 1) "Die Formel, SVERWEIS(), ..." - collapse the split/gramStart runs back
#    into a single run (text itself is unchanged).
Replace-Text "Die Formel, SVERWEIS(), ZÄHLENWENN(), SUMMEWENN()" "Die Formel, SVERWEIS(), ZÄHLENWENN(), SUMMEWENN()"

# 2) "In den Zellen B2:O2 bzw. in den Zellen O2:O4 ..." - merge runs.
Replace-Text "In den Zellen B2:O2 bzw. in den Zellen O2:O4 verwenden Sie bitte das folgende Format:" "In den Zellen B2:O2 bzw. in den Zellen O2:O4 verwenden Sie bitte das folgende Format:"

# 3) "Nutzen Sie dazu die Funktionen ZÄHLENWENN() und SUMMEWENN()" - merge runs.
Replace-Text "Nutzen Sie dazu die Funktionen ZÄHLENWENN() und SUMMEWENN()" "Nutzen Sie dazu die Funktionen ZÄHLENWENN() und SUMMEWENN()"

# 4) "P4: Nutzen Sie die Funktion SUMMENPRODUKT() auf geeignete Weise." - merge runs.
Replace-Text "P4: Nutzen Sie die Funktion SUMMENPRODUKT() auf geeignete Weise." "P4: Nutzen Sie die Funktion SUMMENPRODUKT() auf geeignete Weise."

# 5) Heading "Klausur" -> "Verkäufe" (whole word only - "Klausurteilen"
#    a few lines down must stay untouched).
Replace-WholeWordText "Klausur" "Verkäufe"

# 6) Shorten the paragraph after the heading: drop the trailing sentences.
Replace-Text "Weiter rechts im Arbeitsblatt sehen Sie für verschiedene Studierende Teilnoten aus zwei Klausurteilen bzw. dem Praktikum. In allen drei Tabellen sind die gleichen Studierenden gegeben. Wichtig: Alle Formeln in den Spalten B:H bzw. J und L müssen automatisch ausfüllbar sein!" ""

# 7) "In den Zellen B2:U2 bzw. in den Zellen W2:AA2 ..." -> new cell ranges.
Replace-Text "In den Zellen B2:U2 bzw. in den Zellen W2:AA2 verwenden Sie bitte das folgende Format:" "In den Zellbereichen A2:G2, I2:M2, O2:S2, U2 und W2:AA2 verwenden Sie bitte das folgende Format:"

# 8) "In der Spalte A passen Sie ..." -> insert "ab A3 ".
Replace-Text "In der Spalte A passen Sie das Zahlenformat so an, dass dort ""Ver-0001"" o.ä. steht." "In der Spalte A ab A3 passen Sie das Zahlenformat so an, dass dort ""Ver-0001"" o.ä. steht."

# 9) Append a clarifying sentence to the "Mitarbeiternummer" bullet.
Replace-Text "Tragen Sie eine Mitarbeiternummer in die Zelle V2 ein." "Tragen Sie eine Mitarbeiternummer in die Zelle V2 ein. (z.B. 3, aber gleiches Format wie in Spalte B."

# 10) The SUMMEWENNS/ZÄHLENWENNS bullet: text content is unchanged (only the
#     lastRenderedPageBreak cache hint shifts internally, which Word
#     recomputes during repagination rather than something COM automation
#     sets), so just collapse the split/gramStart runs.
Replace-Text "In den drei Tabellen ""Verkäufe pro Produkt und Verkaufsgebiet"", ""Umsatz pro Produkt und Verkaufsgebiet"" sowie ""Umsatz pro Produkt und Verkaufsgebiet und Mitarbeiter"" nutzen Sie die Formeln SUMMEWENNS() bzw. ZÄHLENWENNS() um die gewünschten Werte zu ermitteln. Passen Sie die Zelle V2 an, und kontrollieren Sie, ob sich die Werte in der dritten Tabelle ändern." "In den drei Tabellen ""Verkäufe pro Produkt und Verkaufsgebiet"", ""Umsatz pro Produkt und Verkaufsgebiet"" sowie ""Umsatz pro Produkt und Verkaufsgebiet und Mitarbeiter"" nutzen Sie die Formeln SUMMEWENNS() bzw. ZÄHLENWENNS() um die gewünschten Werte zu ermitteln. Passen Sie die Zelle V2 an, und kontrollieren Sie, ob sich die Werte in der dritten Tabelle ändern."

Write-Output "done"
